# Tambah status pengiriman & tampilan tabel berwarna
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that *looks* numeric (long NIK codes, phone numbers,
# the "Jumlah" amount) but must stay stored as text, same as the rest of the
# sheet's NIK/No_HP/Jumlah columns.
function Set-TextCell {
    param($row, $col, $value)
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

function Set-Cell {
    param($row, $col, $value)
    $ws.Cells.Item($row, $col).Value = $value
}

# --- Row 5: Jumlah (E5) was stored as a number (65000), now stored as text "65000" ---
Set-TextCell 5 5 "65000"

# --- Row 6 (new) ---
Set-TextCell 6 1  "1234456278949542"
Set-Cell     6 2  "BG4576HI"
Set-Cell     6 3  "Nia Rahmadani"
Set-Cell     6 4  "02-08-2025 16:59"
Set-TextCell 6 5  "65000"
Set-Cell     6 6  "Bank Mandiri"
Set-Cell     6 7  "Nia Rahmadani"
Set-TextCell 6 8  "085267947261"
Set-Cell     6 9  "Palembang"
Set-Cell     6 10 "J&T"

# --- Row 7 (new) ---
Set-TextCell 7 1  "1234456278949542"
Set-Cell     7 2  "BG4576HI"
Set-Cell     7 3  "Nia Rahmadani"
Set-Cell     7 4  "02-08-2025 17:13"
Set-TextCell 7 5  "65000"
Set-Cell     7 6  "Bank Rakyat Indonesia (BRI)"
Set-Cell     7 7  "Nia Rahmadani"
Set-TextCell 7 8  "085267947261"
Set-Cell     7 9  "Jl. Melati, Palembang"
Set-Cell     7 10 "JNE"

# --- Row 8 (new) - Jumlah (E8) is a genuine number this time ---
Set-TextCell 8 1  "1234456278949542"
Set-Cell     8 2  "BG4576HI"
Set-Cell     8 3  "Nia Rahmadani"
Set-Cell     8 4  "02-08-2025 17:38"
Set-Cell     8 5  65000
Set-Cell     8 6  "Bank Rakyat Indonesia (BRI)"
Set-Cell     8 7  "Nia Rahmadani"
Set-TextCell 8 8  "085267947261"
Set-Cell     8 9  "Jl. Melati, Palembang"
Set-Cell     8 10 "JNE"
